# Apply Malay translations to the "Extra one page for Facilitator Manual" doc
$d = $word.ActiveDocument

function Replace-InParagraph($index, $old, $new) {
    $r = $d.Paragraphs.Item($index).Range
    $r.Find.Execute($old, $true, $false, $false, $false, $false, `
                     $true, 1, $false, $new, 2)
}

Replace-InParagraph 4  "Stages of Child Development " "Peringkat Perkembangan Anak "
Replace-InParagraph 5  "Toddler (2-3 years old) " "Kanak-kanak bertatih (2-3 tahun) "
Replace-InParagraph 6  "Food, sleep, cleaning, comfort and safety." "Makanan, tidur, kebersihan, keselesaan dan keselamatan."
Replace-InParagraph 7  "Strong attachment (bonding) with parent/ caregiver." "Hubungan erat (ikatan) dengan ibu bapa/ penjaga."
Replace-InParagraph 8  "Stimulation and attention." "Rangsangan dan perhatian."
Replace-InParagraph 9  "Becomes more curious; wants to explore and become more independent." "Sikap ingin tahu yang tinggi; suka meneroka dan menjadi lebih berdikari."
Replace-InParagraph 10 "Wants to learn how to do new things (e.g. dress and undress) and wants to make own decisions." "Suka belajar melakukan perkara baharu (cth memakai dan menanggalkan pakaian) dan ingin membuat keputusan sendiri."
Replace-InParagraph 39 "May be very active." "Mungkin jadi sangat aktif."
